$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.24"
$ws.Range("E2").Value = "'6.19%"
$ws.Range("D3").Value = "'40.15"
$ws.Range("E3").Value = "'7.36%"
$ws.Range("D4").Value = "'5.587"
$ws.Range("E4").Value = "'9.05%"
$ws.Range("D5").Value = "'0.08114"
$ws.Range("E5").Value = "'3.44%"
$ws.Range("D6").Value = "'4.552"
$ws.Range("E6").Value = "'3.62%"
$ws.Range("D7").Value = "'8.672"
$ws.Range("E7").Value = "'4.93%"
$ws.Range("D8").Value = "'1.991"
$ws.Range("E8").Value = "'6.12%"
$ws.Range("E9").Value = "'1.17%"
$ws.Range("D10").Value = "'0.9498"
$ws.Range("E10").Value = "'2.80%"
$ws.Range("D11").Value = "'0.1277"
$ws.Range("E11").Value = "'11.58%"
$ws.Range("D12").Value = "'0.1970"
$ws.Range("E12").Value = "'3.35%"
$ws.Range("D13").Value = "'0.09199"
$ws.Range("E13").Value = "'4.00%"
$ws.Range("E14").Value = "'7.45%"
$ws.Range("D15").Value = "'0.09591"
$ws.Range("E15").Value = "'-0.06%"
$ws.Range("D16").Value = "'0.001326"
$ws.Range("E16").Value = "'-3.86%"
$ws.Range("D17").Value = "'0.006261"
$ws.Range("E17").Value = "'4.41%"
$ws.Range("E18").Value = "'-0.77%"
$ws.Range("E19").Value = "'1.48%"
$ws.Range("D20").Value = "'7.473"
$ws.Range("E20").Value = "'18.00%"
$ws.Range("D21").Value = "'0.1352"
$ws.Range("E21").Value = "'4.27%"
$ws.Range("D23").Value = "'0.04437"
$ws.Range("E23").Value = "'1.87%"
$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'2.39%"
$ws.Range("D25").Value = "'0.004276"
$ws.Range("E25").Value = "'-0.11%"
$ws.Range("D26").Value = "'0.0001202"
$ws.Range("E26").Value = "'-14.21%"
$ws.Range("D27").Value = "'0.0003996"
$ws.Range("E27").Value = "'37.87%"
$ws.Range("D39").Value = "'0.02531"
$ws.Range("E39").Value = "'16.94%"
$ws.Range("D40").Value = "'0.05207"
$ws.Range("E40").Value = "'3.97%"
$ws.Range("D41").Value = "'0.007730"
$ws.Range("E41").Value = "'1.98%"
$ws.Range("E42").Value = "'5.68%"
$ws.Range("D43").Value = "'0.008896"
$ws.Range("E43").Value = "'4.56%"
$ws.Range("D44").Value = "'0.002193"
$ws.Range("E44").Value = "'9.26%"
$ws.Range("D45").Value = "'0.009640"
$ws.Range("E45").Value = "'18.57%"
$ws.Range("D46").Value = "'0.00006667"
$ws.Range("E46").Value = "'2.15%"
$ws.Range("E47").Value = "'0.10%"
$ws.Range("D48").Value = "'0.002881"
$ws.Range("E48").Value = "'-12.56%"
$ws.Range("D49").Value = "'0.002304"
$ws.Range("E49").Value = "'59.74%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.10%"
